# Reading test case steps method
# Fill in the "Test Case Name" step value for the first step row of TC001
# and make TC001 the active sheet/selection, matching how the workbook
# was left after running the "reading test case steps" routine.

$wb = $excel.ActiveWorkbook

$tc001 = $wb.Worksheets.Item("TC001")

# Record that step's test case name ("Log In as Admin") also applies to
# the row describing the second step (row 3) of TC001.
$tc001.Range("B3").Value = "Log In as Admin"

# Leave TC001 as the active/selected sheet with B3 selected, as it was
# after the edit (Matrix sheet was active/tabSelected before this change).
$tc001.Activate()
$tc001.Range("B3").Select()
